$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-04 Wednesday", "2024-12-05 Thursday"),
    @("765÷7=", "196÷5="),
    @("634÷5=", "923÷3="),
    @("686÷7=", "670÷4="),
    @("317÷6=", "115÷3="),
    @("923÷5=", "325÷9="),
    @("543÷5=", "507÷6="),
    @("956÷6=", "227÷2="),
    @("973÷7=", "887÷5="),
    @("577÷2=", "521÷6="),
    @("182÷9=", "176÷6="),
    @("934÷8=", "169÷3="),
    @("919÷8=", "725÷3="),
    @("940÷7=", "876÷9="),
    @("678÷4=", "814÷7="),
    @("931÷3=", "755÷8="),
    @("219÷5=", "325÷2="),
    @("140÷8=", "611÷5="),
    @("694÷3=", "311÷9="),
    @("987÷7=", "619÷7="),
    @("901÷9=", "613÷3="),
    @("421÷2=", "897÷3="),
    @("877÷4=", "978÷4="),
    @("497÷5=", "418÷2="),
    @("933÷6=", "778÷5="),
    @("613÷2=", "166÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
